# Project configuration workbook update:
# - Replace "individualParamsFile" / "IndividualParameters.xlsx" row with
#   "individualsFile" / "Individuals.xlsx" (the individual-specific parameters
#   are now part of the "Individuals" file)
# - Remove the now-obsolete "individualPhysiologyFile" / "IndividualBiometrics.xlsx" row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5: individualParamsFile -> individualsFile
$ws.Range("A5").Value = "individualsFile"
$ws.Range("B5").Value = "Individuals.xlsx"
# Description (C5) stays the same - "Name of the excel file with individual-specific
# model parametrization. Must be located in the "paramsFolder""

# Remove row 6 (individualPhysiologyFile / IndividualBiometrics.xlsx), shifting
# all following rows up by one
$ws.Rows.Item(6).Delete()

# Leave selection on A8, matching where the cursor ended up after the edit
$ws.Range("A8").Select()
